# "Fuel" sheet update for SIN database
# - updated PEN & CO2 for natural gas (NG) in "Fuel" sheet for SIN
#   (data source: ecoinvent 3.4 market for natural gas, burned in gas motor, for storage, GLO)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")

# PEN (Primary Energy Need) for Natural Gas (row 2, column C):
# previously a hard-coded value of 1.403, now computed from ecoinvent 3.4 factors.
$ws.Range("C2").Formula = "=1.1767+0.0019487+0.0000015726"

# CO2 for Natural Gas (row 2, column D): updated value from ecoinvent 3.4
$ws.Range("D2").Value = 0.06682

# Update the "reference" cell (row 2, column F) to point to the new data source description
$ws.Range("F2").Value = "ecoinvent 3.4 - market for natural gas, burned in gas motor, for storage_GLO_2017_Allocation, cut-off"

# Reflect the active cell the author was last working on in this sheet
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
